$d = $word.ActiveDocument

# 1) Delete the entire "potatoes" list paragraph (text + its paragraph mark).
#    This lets "bread crumbs even" inherit the vacated list slot/indent level.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq "potatoes") {
        $p.Range.Delete()
        break
    }
}

# 2) Word's auto-managed "_GoBack" bookmark follows the location of the last
#    edit. Relocate it from its old spot to the start of the now-adjacent
#    "bread crumbs even" paragraph (which is where the edit above took place).
$d.Bookmarks.Item("_GoBack").Delete()
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq "bread crumbs even") {
        $target = $d.Range($p.Range.Start, $p.Range.Start)
        $d.Bookmarks.Add("_GoBack", $target)
        break
    }
}
